# Updates cryptos price (D) and volume-change (E) columns to match the
# refreshed source data for this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.252.76"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "'1.566.65"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'210.78"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "'22.10"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "'0.0598"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "'1.794.16"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "'1.576.42"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").Value = "'3.78"
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "'0.519"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "'27.213.89"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "'62.17"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "'7.51"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "'217.31"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'0.0₃0702"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "'9.22"
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "'153.36"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("D26").Value = "'6.64"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'15.06"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("E28").Value = "  +1.92%  "
$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +2.32%  "
$ws.Range("D31").Value = "'0.0472"
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("D34").Value = "'1.443.60"
$ws.Range("E35").Value = "  +3.55%  "
$ws.Range("D36").Value = "'1.61"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'0.534"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").Value = "'5.90"
$ws.Range("E40").Value = "  +2.14%  "
$ws.Range("D41").Value = "'0.810"
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").Value = "'1.01"
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D45").Value = "'64.50"
$ws.Range("E45").Value = "  -0.44%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'1.705.41"
$ws.Range("E47").Value = "  +0.42%  "
$ws.Range("D48").Value = "'86.02"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").Value = "'0.0956"
$ws.Range("E51").Value = "  -0.30%  "
